$wb = $excel.ActiveWorkbook

# DRONE1 keeps its data and becomes the single surviving "DRONES" sheet.
$ws1 = $wb.Worksheets.Item("DRONE1")

# Insert a new first column to hold the drone's Name.
$ws1.Range("A1:A2").EntireColumn.Insert()
$ws1.Range("A1:A4").ColumnWidth = 13.830729166666666

$ws1.Range("A1").Value = "Name"
$ws1.Range("A2").Value = "DJI Phantom 4"

$ws1.Range("A3").Value = "DJI Mavic 3"
$ws1.Range("B3").Value = 25.6
$ws1.Range("C3").Value = 0.15
$ws1.Range("D3").Value = 0.98
$ws1.Range("E3").Value = 0.54

$ws1.Range("A4").Value = " Custom Drone 1"
$ws1.Range("B4").Value = 5.02
$ws1.Range("C4").Value = 0.23
$ws1.Range("D4").Value = 0.17
$ws1.Range("E4").Value = 0.58

# Drop the other two sheets; re-fetch by name right before deleting so the
# reference isn't stale after the previous delete shifts the collection.
$wb.Worksheets.Item("DRONE2").Delete()
$wb.Worksheets.Item("DRONE3").Delete()

$ws1.Name = "DRONES"
$ws1.Select()
$ws1.Range("E5").Select()
